$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Update the time_taken (F column) timestamps on the existing "data" sheet ---
$dataSheet.Range("F2").Value = "2021-10-05 14:19:28.306965"
$dataSheet.Range("F3").Value = "2021-10-05 14:19:28.306973"
$dataSheet.Range("F4").Value = "2021-10-05 14:19:28.306977"
$dataSheet.Range("F5").Value = "2021-10-05 14:19:28.306979"

# --- Add the new "metadata" worksheet after "data" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$metaSheet = $wb.Worksheets.Add($null, $lastSheet)
$metaSheet.Name = "metadata"

# Reuse the header style (bold, centered, bordered) from the "data" sheet's header row
$dataSheet.Range("B1:F1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)

# Reuse the index-column style from the "data" sheet
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)

# Header row
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Cerebral folate deficiency"
$metaSheet.Range("C2").Value = 109

# D2 needs to be stored as text "1.2" (not a number) with no special number format,
# so build it as a text formula first, then paste-special values-only over itself.
$metaSheet.Range("D2").Formula = "=""1.2"""
$metaSheet.Range("D2").Copy()
$metaSheet.Range("D2").PasteSpecial(-4163)

$metaSheet.Range("E2").Value = "2017-11-05T02:37:20.047324Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:19:28.303270"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/109/?format=json"

$dataSheet.Select()
